$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New product rows appended after the existing data (rows 35-41).
# Columns: A=kind, B=ref_code, C=sku, D=condition, E=warranty, F=sale_price, G=stock

$rows = @(
    @{ Row=35; Kind="processor"; RefCode="i3-9100";  Sku="i3-9100/tray";  Condition="Tray"; Warranty="3m"; Price=2200000; Stock=10 },
    @{ Row=36; Kind="processor"; RefCode="i3-9100T"; Sku="i3-9100t/tray"; Condition="Tray"; Warranty="3m"; Price=3000000; Stock=10 },
    @{ Row=37; Kind="processor"; RefCode="i3-9100F"; Sku="i3-9100f/tray"; Condition="Tray"; Warranty="3m"; Price=2000000; Stock=10 },
    @{ Row=38; Kind="barebone";  RefCode="dell-optiplex-7060-sff"; Sku="dell-optiplex-7060-sff/used"; Condition="Used"; Warranty="3m"; Price=1850000; Stock=10 },
    @{ Row=39; Kind="gpu";       RefCode="nvidia-quadro-k600"; Sku="nvidia-quadro-k600/used"; Condition="Used"; Warranty="3m"; Price=600000; Stock=10 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.Kind
    $ws.Range("B$n").Value = $r.RefCode
    $ws.Range("C$n").Value = $r.Sku
    $ws.Range("D$n").Value = $r.Condition
    $ws.Range("E$n").Value = $r.Warranty
    $ws.Range("F$n").Value = $r.Price
    $ws.Range("G$n").Value = $r.Stock
}

# The last two memory rows were upserted in bulk: both ref_codes were
# written first, then both computed skus - so the shared-string table
# picks up ref_codes before skus for this pair.
$ws.Range("A40").Value = "memory"
$ws.Range("A41").Value = "memory"
$ws.Range("B40").Value = "4gb-dimm-ddr4-2666-mixed"
$ws.Range("B41").Value = "8gb-dimm-ddr4-2666-mixed"
$ws.Range("C40").Value = "4gb-dimm-ddr4-2666-mixed/used"
$ws.Range("C41").Value = "8gb-dimm-ddr4-2666-mixed/used"
$ws.Range("D40").Value = "Used"
$ws.Range("D41").Value = "Used"
$ws.Range("E40").Value = "3m"
$ws.Range("E41").Value = "3m"
$ws.Range("F40").Value = 340000
$ws.Range("F41").Value = 680000
$ws.Range("G40").Value = 10
$ws.Range("G41").Value = 10

# Narrow the stock column now that it only needs to fit small numbers.
$ws.Columns.Item(7).ColumnWidth = 5.83

# Move the active selection to the last-edited row, like a user who just
# finished typing the new data.
$ws.Range("C40").Select() | Out-Null
